# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Periodo Mora" values from 2507 -> 2508 for existing worker rows ---
$ws.Range("E16").Value = "2508"
$ws.Range("E17").Value = "2508"
$ws.Range("E18").Value = "2508"

# --- Remove the worker row for ANDERSON DUVAN REALES OSORNO (row 17) ---
# (this shifts the SANDRA MILENA MEDELLIN ARIZA row and the signature rows up by one)
$ws.Rows.Item(17).Delete()

# --- Update the account summary figures to reflect the removed worker ---
$ws.Range("E11").Value = 113880
$ws.Range("C13").Value = 2

# --- Narrow column D now that the longest name it needs to fit is shorter ---
$ws.Columns.Item(4).ColumnWidth = 30
